$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export_Contract")
$ws.Rows.Item(22).Delete()
